$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark from the title paragraph.
# ------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) "... for HCl and for CO." -> "... (in GHz) for HCl and for CO."
# ------------------------------------------------------------------
$d.Content.Find.Execute(" for HCl and for CO.", $false, $false, $false, $false, $false, $false, 1, $false, " (in GHz) for HCl and for CO.", 2)

# ------------------------------------------------------------------
# 3) "Calculate the reduced mass of the different molecules..."
#    -> "Calculate the reduced mass (in atomic mass units u) of the
#        different molecules..."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Calculate the reduced mass of the different molecules from the masses of the individual atoms.", $false, $false, $false, $false, $false, $false, 1, $false, "Calculate the reduced mass (in atomic mass units u) of the different molecules from the masses of the individual atoms.", 2)

# ------------------------------------------------------------------
# 4) "Calculate the bond length of the various molecules (except O"
#    -> "Calculate the bond length (in pm) of the various molecules
#        (except O" with a "_GoBack" bookmark re-inserted right after
#        "(in pm".
# ------------------------------------------------------------------
$d.Content.Find.Execute("Calculate the bond length of the various molecules (except O", $false, $false, $false, $false, $false, $false, 1, $false, "Calculate the bond length (in pm) of the various molecules (except O", 2)

$r = $d.Content
$r.Find.Execute("(in pm", $false, $false, $false, $false, $false, $false, 1, $false, "", 0)
$pos = $r.End
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
